$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-29 Wednesday" "2025-01-30 Thursday"

Replace-Text "53÷7=" "79÷3="
Replace-Text "25÷8=" "25÷9="
Replace-Text "29÷6=" "11÷7="
Replace-Text "28÷7=" "95÷9="
Replace-Text "58÷7=" "28÷9="

Replace-Text "45÷6=" "63÷2="
Replace-Text "98÷2=" "33÷8="
Replace-Text "93÷2=" "30÷6="
Replace-Text "20÷6=" "60÷7="
Replace-Text "47÷7=" "63÷8="

Replace-Text "22÷4=" "16÷6="
Replace-Text "70÷2=" "21÷4="
Replace-Text "94÷9=" "62÷7="
Replace-Text "90÷6=" "24÷8="
Replace-Text "59÷4=" "33÷7="

Replace-Text "65÷5=" "50÷6="
Replace-Text "56÷7=" "46÷2="
Replace-Text "61÷5=" "63÷8="
Replace-Text "82÷7=" "18÷8="
Replace-Text "99÷5=" "76÷6="

Replace-Text "36÷4=" "77÷6="
Replace-Text "48÷7=" "17÷9="
Replace-Text "19÷9=" "24÷7="
Replace-Text "40÷3=" "59÷3="
Replace-Text "80÷5=" "70÷5="
